$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column (price) cells are forced to text so numeric-looking
# strings such as "60.304.09" or "1.00" are not reinterpreted as numbers/dates.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.304.09"
$ws.Range("E2").Value = "  -5.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.296.80"
$ws.Range("E3").Value = "  -5.18%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.22"
$ws.Range("E5").Value = "  -4.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.00"
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.296.64"
$ws.Range("E8").Value = "  -5.19%  "
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.35"
$ws.Range("E10").Value = "  -4.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.116"
$ws.Range("E11").Value = "  -5.09%  "
$ws.Range("E12").Value = "  -4.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.860.68"
$ws.Range("E13").Value = "  -5.09%  "
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.301.40"
$ws.Range("E15").Value = "  -5.06%  "
$ws.Range("E16").Value = "  -5.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.489.87"
$ws.Range("E17").Value = "  -5.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.14"
$ws.Range("E18").Value = "  -3.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.65"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.30"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.94"
$ws.Range("E21").Value = "  -10.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "350.30"
$ws.Range("E22").Value = "  -9.12%  "
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.428.20"
$ws.Range("E25").Value = "  -5.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.17"
$ws.Range("E26").Value = "  -7.26%  "
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("E29").Value = "  +2.41%  "
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.81"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.151"
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.10"
$ws.Range("E33").Value = "  -5.82%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.324.85"
$ws.Range("E35").Value = "  -5.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.59"
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.25"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.75"
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.47"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "157.54"
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("E41").Value = "  -3.88%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.91"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.32"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.741"
$ws.Range("E45").Value = "  -7.08%  "
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("E48").Value = "  -4.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.66"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.57"
$ws.Range("E50").Value = "  +4.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.858"
$ws.Range("E51").Value = "  -5.03%  "
